# ---------------------------------------------------------------------------
# Rebuilds the QA data workbook:
#   1. Duplicates the original "QA_Worksheet1" sheet as a new trailing sheet
#      named "Sheet1" (the old/simple 4-row test data), pointing its two
#      hyperlink cells at the refreshed support-article URLs.
#   2. Repopulates "QA_Worksheet1" itself with the expanded 13-row article
#      list (Trying to Join / During Your Webinar categories), with a single
#      live hyperlink on C5.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item(1)

# --- 1. Duplicate the current sheet to the end, rename it "Sheet1" ----------
$orig.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dup = $wb.Worksheets.Item($wb.Worksheets.Count)
$dup.Name = "Sheet1"

# Drop the old hyperlinks on the duplicate and point the cells at the new
# article URLs (keep their existing "hyperlink-look" style).
$dup.Range("C3").Hyperlinks.Delete()
$dup.Range("C4").Hyperlinks.Delete()
$dup.Range("C3").Value = "http://support.citrixonline.com/en_US/Webinar/all_files/G2W090002?__col_mkt_cookies=__col_visit"
$dup.Range("C4").Value = "http://support.citrixonline.com/en_US/Webinar/all_files/G2W030004?__col_mkt_cookies=__col_visit"
$dup.Range("D2").Select()

# --- 2. Rebuild QA_Worksheet1 with the expanded article list ---------------
$ws = $orig

# Remove the two existing hyperlinks (URLs/text are all changing anyway) and
# drop the "hyperlink look" styling those two cells carried.
$ws.Range("C3").Hyperlinks.Delete()
$ws.Range("C4").Hyperlinks.Delete()
$ws.Range("C3:C4").Style = "Normal"

# Apply the "category" column formatting (style of B1) down through the new
# rows before writing values, so every B cell in 2:13 matches the original
# look.
$ws.Range("B1").Copy()
$ws.Range("B2:B13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$rows = @(
  @("Join Page - Trying to Join - Article 1 Test", "Trying to Join", "http://support.citrixonline.com/en_US/Webinar/all_files/G2W060002?__col_mkt_cookies=__col_visit", "What are some quick fixes I can try for joining?"),
  @("Join Page - Trying to Join - Article 2 Test", "Trying to Join", "http://support.citrixonline.com/en_US/Webinar/all_files/G2W090002?__col_mkt_cookies=__col_visit", "A firewall is preventing me from joining"),
  @("Join Page - Trying to Join - Article 3 Test", "Trying to Join", "http://support.citrixonline.com/en_US/Webinar/all_files/G2W030004?__col_mkt_cookies=__col_visit", "I am having Download issues"),
  @("Join Page - Trying to Join - Article 4 Test", "Trying to Join", "http://support.citrixonline.com/en_US/Webinar/help_files/G2W060005?__col_mkt_cookies=__col_visit", "How do I install GoToWebinar on a Mac?"),
  @("Join Page - Trying to Join - Article 5 Test", "Trying to Join", "http://support.citrixonline.com/en_US/Webinar/help_files/G2W060021?__col_mkt_cookies=__col_visit", "How do I install GoToWebinar on a PC?"),
  @("Join Page - Trying to Join - Article 6 Test", "Trying to Join", "http://support.citrixonline.com/en_US/webinar/all_files/G2W060007?__col_mkt_cookies=__col_visit", "I'm still having trouble joining - what else can I try?"),
  @("Join Page - During Your Webinar - Article 1 Test", "During Your Webinar", "https://support.citrixonline.com/en_US/webinar/knowledge_articles/000025306?__col_mkt_cookies=__col_visit", 'I registered for the webinar but now I see "The session is full"'),
  @("Join Page - During Your Webinar - Article 2 Test", "During Your Webinar", "http://support.citrixonline.com/en_US/webinar/help_files/G2W090003?__col_mkt_cookies=__col_visit", "I can only hear the audio or see the screen, but not both"),
  @("Join Page - During Your Webinar - Article 3 Test", "During Your Webinar", "http://support.citrixonline.com/en_US/webinar/knowledge_articles/000161527?title=Waiting+for+Organizer&__col_mkt_cookies=__col_visit", "It says I've connected, but the webinar isn't started"),
  @("Join Page - During Your Webinar - Article 4 Test", "During Your Webinar", "http://support.citrixonline.com/en_US/Webinar/all_files/G2W050039?__col_mkt_cookies=__col_visit", "How do I unmute myself?"),
  @("Join Page - During Your Webinar - Article 5 Test", "During Your Webinar", "http://support.citrixonline.com/en_US/Webinar/all_files/G2W050053?__col_mkt_cookies=__col_visit", "Why can't anyone hear me?"),
  @("Join Page - During Your Webinar - Article 6 Test", "During Your Webinar", "https://support.citrixonline.com/en_US/webinar/knowledge_articles/000064729?__col_mkt_cookies=__col_visit", "I have a question or feedback about the presentation itself (materials, certification, etc.)")
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

# Live hyperlink only on C5 (first "Trying to Join" article-4 row) - added
# without the usual blue/underline "hyperlink style".
$ws.Hyperlinks.Add($ws.Range("C5"), "http://support.citrixonline.com/en_US/Webinar/help_files/G2W060005?__col_mkt_cookies=__col_visit")
$ws.Range("C5").Style = "Normal"

# Column widths matching the refreshed content (A manually sized; B-D sized
# to fit the new, longer article names/urls/questions).
$ws.Columns.Item(1).ColumnWidth = 40.666666666666664
$ws.Columns.Item(2).ColumnWidth = 17.833333333333332
$ws.Columns.Item(3).ColumnWidth = 116.16666666666667
$ws.Columns.Item(4).ColumnWidth = 50.333333333333336

$ws.Range("B17").Select()
